$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank separator rows (14, 25, 35, 40) so the data becomes
# contiguous again. Deleting from the bottom up keeps the remaining row
# numbers stable while each delete is applied.
$ws.Rows(40).Delete()
$ws.Rows(35).Delete()
$ws.Rows(25).Delete()
$ws.Rows(14).Delete()

# Reset the view: scroll back to the top and select cell B1.
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("B1").Select()
